# Starlendar.xlsx: rename the weekday "Sumday" to "Zingday".
#
# While retyping the six weekday header labels (row 7-12, column A) on both
# the "Decimal" and "Dozenal" sheets, the author also picked up a stray
# leading space on each of them, and replaced "Sumday" with " Zingday".
# The "Ultra days AKA intercalary days" header on the Decimal sheet got
# shortened to match the Dozenal sheet's existing "Ultra days" header.
# Finally the workbook was left with the Dozenal tab active/selected.

$wb = $excel.ActiveWorkbook

# ---- Decimal sheet ----
$ws1 = $wb.Worksheets.Item("Decimal")
$ws1.Range("A7").Value  = " Starday"
$ws1.Range("A8").Value  = " Topday"
$ws1.Range("A9").Value  = " Aceday"
$ws1.Range("A10").Value = " Rolday"
$ws1.Range("A11").Value = " Zingday"
$ws1.Range("A12").Value = " Funday"
$ws1.Range("R1").Value  = "Ultra days"
[void]$ws1.Range("T18").Select()

# ---- Dozenal sheet ----
$ws2 = $wb.Worksheets.Item("Dozenal")
$ws2.Range("A7").Value  = " Starday"
$ws2.Range("A8").Value  = " Topday"
$ws2.Range("A9").Value  = " Aceday"
$ws2.Range("A10").Value = " Rolday"
$ws2.Range("A11").Value = " Zingday"
$ws2.Range("A12").Value = " Funday"

# Leave the workbook with the Dozenal sheet active/selected, matching the
# final saved UI state.
[void]$ws2.Activate()
[void]$ws2.Range("T19").Select()
